$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append, matching the style/formatting of the last existing
# data row (row 18).
$newRows = @(
    @("26/04/2025 14:53:14", "26°", "65%"),
    @("26/04/2025 14:55:15", "26°", "65%"),
    @("26/04/2025 14:55:23", "26°", "65%")
)

$templateRow = 18
$startRow = 19

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Copy the formatting (style) of the template row onto the new row first.
    $ws.Range("A$templateRow`:C$templateRow").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)

    # Force the third column to text so values like "65%" aren't
    # reinterpreted as a percentage number.
    $ws.Cells.Item($r, 3).NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]

    # Re-apply the template row's formatting so the number-format tweak above
    # doesn't leave the new row on a different style index than the template.
    $ws.Range("A$templateRow`:C$templateRow").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)
}
